$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.07851779460907
$ws.Range("B1").Value = 2.5318443775177
$ws.Range("C1").Value = 4.784962177276611
$ws.Range("D1").Value = 2.338709831237793
$ws.Range("E1").Value = 1.075323820114136
